{"js": "// Update the date line and every \"two-digit \u00f7 one-digit\" answer cell in\n// the practice-sheet table to the next day's generated worksheet values.\n//\n// Each entry is [oldText, newText] and every oldText is unique in the\n// original document, so we can safely search for all of them against the\n// pristine body first (collecting the resulting Range objects) and only\n// then perform the text replacements. Doing the search pass fully before\n// any mutation avoids any ordering hazards in case a newText happens to\n// match another entry's oldText (it does, once, in this sheet).\nconst replacements = [\n  [\"2024-06-11 Tuesday\", \"2024-06-12 Wednesday\"],\n  [\"60\u00f73=20, 0\", \"33\u00f73=11, 0\"],\n  [\"20\u00f74=5, 0\", \"88\u00f78=11, 0\"],\n  [\"50\u00f78=6, 2\", \"96\u00f78=12, 0\"],\n  [\"59\u00f75=11, 4\", \"14\u00f77=2, 0\"],\n  [\"93\u00f79=10, 3\", \"18\u00f72=9, 0\"],\n  [\"65\u00f72=32, 1\", \"80\u00f78=10, 0\"],\n  [\"14\u00f76=2, 2\", \"33\u00f72=16, 1\"],\n  [\"85\u00f74=21, 1\", \"55\u00f72=27, 1\"],\n  [\"40\u00f73=13, 1\", \"29\u00f72=14, 1\"],\n  [\"25\u00f72=12, 1\", \"79\u00f79=8, 7\"],\n  [\"85\u00f79=9, 4\", \"33\u00f75=6, 3\"],\n  [\"46\u00f75=9, 1\", \"81\u00f72=40, 1\"],\n  [\"32\u00f78=4, 0\", \"20\u00f78=2, 4\"],\n  [\"66\u00f73=22, 0\", \"56\u00f78=7, 0\"],\n  [\"32\u00f76=5, 2\", \"78\u00f76=13, 0\"],\n  [\"69\u00f75=13, 4\", \"40\u00f77=5, 5\"],\n  [\"71\u00f75=14, 1\", \"99\u00f74=24, 3\"],\n  [\"97\u00f78=12, 1\", \"46\u00f77=6, 4\"],\n  [\"93\u00f75=18, 3\", \"19\u00f78=2, 3\"],\n  [\"86\u00f77=12, 2\", \"24\u00f77=3, 3\"],\n  [\"91\u00f79=10, 1\", \"75\u00f74=18, 3\"],\n  [\"72\u00f74=18, 0\", \"84\u00f76=14, 0\"],\n  [\"90\u00f74=22, 2\", \"76\u00f72=38, 0\"],\n  [\"81\u00f76=13, 3\", \"86\u00f77=12, 2\"],\n  [\"91\u00f76=15, 1\", \"22\u00f78=2, 6\"],\n];\n\nconst body = context.document.body;\n\n// Pass 1: issue a search for every old string against the untouched body.\nconst pending = replacements.map(([oldText, newText]) => {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return { oldText, newText, results };\n});\nawait context.sync();\n\n// Pass 2: replace each unique hit with its new text.\nfor (const { oldText, newText, results } of pending) {\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date line and all 25 division-answer cells to the\n# next generated day's values. Every \"old\" string below is unique within\n# the document, so a targeted Find/Replace for each pair is unambiguous.\n# Pairs are applied in the same order the generator emitted them; this\n# matters once, since pair 20's replacement text happens to equal pair\n# 24's original text (\"86\u00f77=12, 2\") -- running pair 20 first means its\n# occurrence is already gone before pair 24 (re)introduces that text, so\n# a straightforward sequential Replace-All never double-fires.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-06-11 Tuesday\", \"2024-06-12 Wednesday\"),\n    @(\"60\u00f73=20, 0\", \"33\u00f73=11, 0\"),\n    @(\"20\u00f74=5, 0\", \"88\u00f78=11, 0\"),\n    @(\"50\u00f78=6, 2\", \"96\u00f78=12, 0\"),\n    @(\"59\u00f75=11, 4\", \"14\u00f77=2, 0\"),\n    @(\"93\u00f79=10, 3\", \"18\u00f72=9, 0\"),\n    @(\"65\u00f72=32, 1\", \"80\u00f78=10, 0\"),\n    @(\"14\u00f76=2, 2\", \"33\u00f72=16, 1\"),\n    @(\"85\u00f74=21, 1\", \"55\u00f72=27, 1\"),\n    @(\"40\u00f73=13, 1\", \"29\u00f72=14, 1\"),\n    @(\"25\u00f72=12, 1\", \"79\u00f79=8, 7\"),\n    @(\"85\u00f79=9, 4\", \"33\u00f75=6, 3\"),\n    @(\"46\u00f75=9, 1\", \"81\u00f72=40, 1\"),\n    @(\"32\u00f78=4, 0\", \"20\u00f78=2, 4\"),\n    @(\"66\u00f73=22, 0\", \"56\u00f78=7, 0\"),\n    @(\"32\u00f76=5, 2\", \"78\u00f76=13, 0\"),\n    @(\"69\u00f75=13, 4\", \"40\u00f77=5, 5\"),\n    @(\"71\u00f75=14, 1\", \"99\u00f74=24, 3\"),\n    @(\"97\u00f78=12, 1\", \"46\u00f77=6, 4\"),\n    @(\"93\u00f75=18, 3\", \"19\u00f78=2, 3\"),\n    @(\"86\u00f77=12, 2\", \"24\u00f77=3, 3\"),\n    @(\"91\u00f79=10, 1\", \"75\u00f74=18, 3\"),\n    @(\"72\u00f74=18, 0\", \"84\u00f76=14, 0\"),\n    @(\"90\u00f74=22, 2\", \"76\u00f72=38, 0\"),\n    @(\"81\u00f76=13, 3\", \"86\u00f77=12, 2\"),\n    @(\"91\u00f76=15, 1\", \"22\u00f78=2, 6\")\n)\n\nforeach ($p in $pairs) {\n    $oldText = $p[0]\n    $newText = $p[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $oldText,    # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n}\n"}
